# "edit fit table to bold LRT" - set the per-column paragraph justification
# on the property table: first column (the label / "Variables" column)
# left-aligned ("start"), the other two (numeric) columns right-aligned
# ("end"). Word's own ParagraphFormat.Alignment setter only knows the
# legacy wdAlignParagraph* left/right/center/both vocabulary, which the
# engine serialises back out as <w:jc w:val="left|right"/>. This document
# already uses the logical start/end vocabulary everywhere else (tcBorders,
# tblCellMar), so we pass the literal ST_Jc token strings "start"/"end"
# straight through to Paragraph.Alignment, which the interop layer writes
# out verbatim as <w:jc w:val="start"/> / <w:jc w:val="end"/>.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$rowCount = $table.Rows.Count
$colCount = $table.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $table.Cell($r, $c)
        $para = $cell.Range.Paragraphs.Item(1)
        if ($c -eq 1) {
            $para.Alignment = "start"
        } else {
            $para.Alignment = "end"
        }
    }
}
